$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1716
$ws.Range("F4").Value = 1181
$ws.Range("F6").Value = 148
$ws.Range("F7").Value = 1404
$ws.Range("F8").Value = 66
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 102
$ws.Range("F11").Value = 635
$ws.Range("F12").Value = 132
$ws.Range("F13").Value = 105
$ws.Range("F14").Value = 1310
$ws.Range("F15").Value = 471
$ws.Range("F16").Value = 468
$ws.Range("F18").Value = 25
$ws.Range("F19").Value = 702
$ws.Range("F20").Value = 2548
$ws.Range("F21").Value = 516
$ws.Range("F22").Value = 17
$ws.Range("F25").Value = 173
$ws.Range("F27").Value = 104
$ws.Range("F28").Value = 565
$ws.Range("F29").Value = 921
$ws.Range("F31").Value = 79
$ws.Range("F35").Value = 239

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 723
$ws.Range("F5").Value = 616
$ws.Range("F6").Value = 616
$ws.Range("F12").Value = 273
$ws.Range("F15").Value = 362
$ws.Range("F16").Value = 362
$ws.Range("F26").Value = 232
$ws.Range("F27").Value = 222

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1751
$ws.Range("F5").Value = 2278
$ws.Range("F6").Value = 905
$ws.Range("F9").Value = 1126
$ws.Range("F10").Value = 252

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1751
$ws.Range("F3").Value = 2278
$ws.Range("F5").Value = 1716
$ws.Range("F8").Value = 905
$ws.Range("F9").Value = 1126
$ws.Range("F10").Value = 252
$ws.Range("F12").Value = 723
$ws.Range("F13").Value = 1181
$ws.Range("F15").Value = 148
$ws.Range("F16").Value = 1404
$ws.Range("F17").Value = 616
$ws.Range("F18").Value = 66
$ws.Range("F19").Value = 102
$ws.Range("F20").Value = 635
$ws.Range("F21").Value = 132
$ws.Range("F23").Value = 105
$ws.Range("F24").Value = 471
$ws.Range("F26").Value = 468
$ws.Range("F28").Value = 702
$ws.Range("F29").Value = 2548
$ws.Range("F30").Value = 516
$ws.Range("F31").Value = 17
$ws.Range("F33").Value = 273
$ws.Range("F34").Value = 173
$ws.Range("F35").Value = 104
$ws.Range("F37").Value = 565
$ws.Range("F38").Value = 921
$ws.Range("F39").Value = 362
$ws.Range("F42").Value = 79
$ws.Range("F46").Value = 232
$ws.Range("F47").Value = 222
$ws.Range("F51").Value = 239

